$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for all data rows (2-262)
# from serial date 45204 to 45205 (one day later).
$ws.Range("C2:C262").Value = 45205
